# RSTK-9221_Inspection Order.xlsx
# - Reorders the 4 data rows (Inventory Item / Inspection Quantity / Lot Number)
#   on the "Create Insp Order" sheet: row2<->row4 and row3<->row5.
# - Updates the active selection to rows 2:3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Create Insp Order")

# --- Swap row 2 <-> row 4 (Inventory Item, Inspection Quantity, Lot Number) ---
$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2
$g2 = $ws.Range("G2").Value2

$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2
$g4 = $ws.Range("G4").Value2

$ws.Range("A2").Value2 = $a4
$ws.Range("B2").Value2 = $b4
$ws.Range("G2").Value2 = $g4

$ws.Range("A4").Value2 = $a2
$ws.Range("B4").Value2 = $b2
$ws.Range("G4").Value2 = $g2

# --- Swap row 3 <-> row 5 (Inventory Item, Inspection Quantity) ---
$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2

$a5 = $ws.Range("A5").Value2
$b5 = $ws.Range("B5").Value2

$ws.Range("A3").Value2 = $a5
$ws.Range("B3").Value2 = $b5

$ws.Range("A5").Value2 = $a3
$ws.Range("B5").Value2 = $b3

# --- Update the active selection to rows 2:3 (whole rows) ---
$ws.Rows("2:3").Select()
